# Corporate Customer excel file update
# Adds DRAWDOWN.ACCOUNT / PRIN.LIQ.ACCT / INT.LIQ.ACCT columns to the four
# TDR maturity sheets, refreshes the sample CUSTOMER.ID on each, and moves
# the active tab from TDRAmendmentPreMaturity to TDRBackDatedMaturityFCY_IBG.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: TDRIssuanceMaturityLCY_IBG
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("TDRIssuanceMaturityLCY_IBG")

$ws1.Range("A2").Value = 11871078

$ws1.Range("E1").Value = "DRAWDOWN.ACCOUNT"
$ws1.Range("F1").Value = "PRIN.LIQ.ACCT"
$ws1.Range("G1").Value = "INT.LIQ.ACCT"

$ws1.Range("E2").Value = 5000000513
$ws1.Range("F2").Value = 5000000513
$ws1.Range("G2").Value = 5000000513

$ws1.Columns.Item(5).AutoFit() | Out-Null
$ws1.Columns.Item(6).AutoFit() | Out-Null
$ws1.Columns.Item(7).AutoFit() | Out-Null

$ws1.Range("B1").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: TDRIssuanceMaturityFCY_IBG
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("TDRIssuanceMaturityFCY_IBG")

$ws2.Range("A2").Value = 11871117

$ws2.Range("F1").Value = "DRAWDOWN.ACCOUNT"
$ws2.Range("G1").Value = "PRIN.LIQ.ACCT"
$ws2.Range("H1").Value = "INT.LIQ.ACCT"

$ws2.Range("F2").Value = 5000000518
$ws2.Range("F2").NumberFormat = "@"
$ws2.Range("G2").Value = 5000000518
$ws2.Range("H2").Value = 5000000518

$ws2.Columns.Item(6).AutoFit() | Out-Null
$ws2.Columns.Item(7).AutoFit() | Out-Null
$ws2.Columns.Item(8).AutoFit() | Out-Null

$ws2.Range("F14").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: TDRBackDatedMaturityLCY_IBG
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("TDRBackDatedMaturityLCY_IBG")

$ws3.Range("A2").Value = 11871085

$ws3.Range("G1").Value = "DRAWDOWN.ACCOUNT"
$ws3.Range("H1").Value = "PRIN.LIQ.ACCT"
$ws3.Range("I1").Value = "INT.LIQ.ACCT"

$ws3.Range("G2:I2").NumberFormat = "@"
$ws3.Range("G2").Value = "5000000520"
$ws3.Range("H2").Value = "5000000520"
$ws3.Range("I2").Value = "5000000520"

$ws3.Columns.Item(7).AutoFit() | Out-Null
$ws3.Columns.Item(8).AutoFit() | Out-Null
$ws3.Columns.Item(9).AutoFit() | Out-Null

$ws3.Range("G1:I2").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet 4: TDRBackDatedMaturityFCY_IBG
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("TDRBackDatedMaturityFCY_IBG")

$ws4.Range("A2").Value = 11871119

$ws4.Range("G1").Value = "DRAWDOWN.ACCOUNT"
$ws4.Range("H1").Value = "PRIN.LIQ.ACCT"
$ws4.Range("I1").Value = "INT.LIQ.ACCT"

$ws4.Range("G2:I2").NumberFormat = "@"
$ws4.Range("G2").Value = "5000000521"
$ws4.Range("H2").Value = "5000000521"
$ws4.Range("I2").Value = "5000000521"

$ws4.Columns.Item(7).AutoFit() | Out-Null
$ws4.Columns.Item(8).AutoFit() | Out-Null
$ws4.Columns.Item(9).AutoFit() | Out-Null

$ws4.Range("G20").Select() | Out-Null

# ---------------------------------------------------------------------
# Make TDRBackDatedMaturityFCY_IBG the active tab (was
# TDRAmendmentPreMaturity) and move the active-sheet selection there.
# ---------------------------------------------------------------------
$ws4.Activate() | Out-Null
